$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell H1, matching the style/formatting of the other header
# cells (e.g. G1: bold font, border, centered alignment)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Add new data column H2:H6 with value 0 (plain numeric cells, no special style,
# matching columns like F/G in the data rows)
$ws.Range("H2:H6").Value = 0
